# Actualización automática 2025-10-31 08:30:08

$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual    = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento    = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
$wsVentasPorGrupo.Range("K4").Value = 848.76
$wsVentasPorGrupo.Range("P47").Value = 176.95
$wsVentasPorGrupo.Range("K60").Value = "2 de 58"
$wsVentasPorGrupo.Range("P60").Value = "3 de 58"

# --- Sheet "VENTA MENSUAL" ---
$wsVentaMensual.Range("F4").Value = 848.76
$wsVentaMensual.Range("F47").Value = 3130.46
$wsVentaMensual.Range("F60").Value = 49978.38

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento.Range("D8").Value = 3404.63
$wsCumplimiento.Range("E8").Value = -2888.508126452166
$wsCumplimiento.Range("F8").Value = 6.596562119323665

$wsCumplimiento.Range("D10").Value = 1102.56
$wsCumplimiento.Range("E10").Value = -714.4520164656079
$wsCumplimiento.Range("F10").Value = 2.840858850568574

$wsCumplimiento.Range("D14").Value = 49978.38
$wsCumplimiento.Range("E14").Value = 4505.387749468964
$wsCumplimiento.Range("F14").Value = 0.9173077058439506
